$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header/data cells to reflect the new week of job-search entries
# Row 1
$ws.Range("A1").Value = 'Date'
$ws.Range("B1").Value = 'Position'
$ws.Range("C1").Value = 'Payrate'
$ws.Range("D1").Value = 'Employer name/address/phone/URL'
$ws.Range("E1").Value = 'Person Contacted'
$ws.Range("F1").Value = 'HOW CONTACTED - Web, phone, mail, job fair, networking, etc.'
$ws.Range("G1").Value = 'Results'
# Row 2
$ws.Range("A2").Value = '<2019-03-25 Mon 19:05>'
$ws.Range("B2").Value = 'SQL Server DBA'
$ws.Range("C2").Value = '$90K'
$ws.Range("D2").Value = 'MassHire Framingham Career Centers'
$ws.Range("E2").Value = 'Lorraine Pocon lpocon@masshiremsw.com'
$ws.Range("F2").Value = 'email'
$ws.Range("G2").Value = 'Signed Up for BSR Notifications'
# Row 3
$ws.Range("A3").Value = '<2019-03-25 Mon 19:05>'
$ws.Range("B3").Value = 'SQL Server DBA'
$ws.Range("C3").Value = '$90K'
$ws.Range("D3").Value = 'MassHire Framingham Career Centers'
$ws.Range("E3").Value = 'Ed Lawrence elawrence@masshiremsw.com'
$ws.Range("F3").Value = 'email'
$ws.Range("G3").Value = 'Confirmation of tomorrows meeting'
# Row 4
$ws.Range("A4").Value = '<2019-03-25 Mon 13:50>'
$ws.Range("B4").Value = 'Oracle E-Business Suite SQA Tester'
$ws.Range("D4").Value = 'Dimensional Thinking, LLC'
$ws.Range("E4").Value = 'Amy Borkar via bullhornmail.com'
$ws.Range("F4").Value = 'email'
$ws.Range("G4").Value = 'Follow up requested'
# Row 5
$ws.Range("A5").Value = '<2019-03-25 Mon 16:47>'
$ws.Range("B5").Value = 'VMware Engineer'
$ws.Range("D5").Value = 'NTT DATA Services'
$ws.Range("E5").Value = 'Siddharth Mishra <siddharth.mishra@nttdata.com>'
$ws.Range("F5").Value = 'email'
$ws.Range("G5").Value = 'Follow up requested'
# Row 6
$ws.Range("A6").Value = '<2019-03-25 Mon 17:48>'
$ws.Range("B6").Value = 'QA Engineer'
$ws.Range("D6").Value = 'Axelon Services Corporation'
$ws.Range("E6").Value = 'Kalyani Munamarthi <kalyani.munamarthi@axelon.com>'
$ws.Range("F6").Value = 'email'
$ws.Range("G6").Value = 'Follow up requested'
# Row 7
$ws.Range("A7").Value = '<2019-03-26 Tue 17:57>'
$ws.Range("B7").Value = 'MS Access developer'
$ws.Range("D7").Value = 'ApTask'
$ws.Range("E7").Value = 'Aravind Vennam <vennam@aptask.com>'
$ws.Range("F7").Value = 'email'
$ws.Range("G7").Value = 'First follow-up email'
# Row 8
$ws.Range("A8").Value = '<2019-03-19 Tue 16:31>'
$ws.Range("B8").Value = 'Programmer - AllStar Staffing Group - Boston, MA'
$ws.Range("D8").Value = 'AllStar Staffing Group'
$ws.Range("E8").Value = 'Rob Greco'
$ws.Range("F8").Value = '215-944-8140'
$ws.Range("G8").Value = 'Talked to Rob and email new resume'
# Row 9
$ws.Range("A9").Value = '<2019-03-26 Tue 18:20>'
$ws.Range("B9").Value = 'DB2 DBA'
$ws.Range("D9").Value = 'Voto Consulting LLC'
$ws.Range("E9").Value = 'Himanshu Kumar himanshu@votoconsulting.com'
$ws.Range("F9").Value = 'email'
$ws.Range("G9").Value = 'Follow up requested'
# Row 10
$ws.Range("A10").Value = '<2019-03-26 Tue 18:26>'
$ws.Range("B10").Value = 'Help Desk Technical Support'
$ws.Range("D10").Value = 'http://www.sigmainc.com'
$ws.Range("E10").Value = 'Surya Ponnada sponnada@sigmainc.com'
$ws.Range("F10").Value = 'email'
$ws.Range("G10").Value = 'Follow up requested'
# Row 11
$ws.Range("A11").Value = '<2019-03-26 Tue 18:30>'
$ws.Range("B11").Value = 'Software Developer (2) C#/.net'
$ws.Range("D11").Value = 'Robert Half Technology'
$ws.Range("E11").Value = 'Ryan Staab ryan.staab@rht.com'
$ws.Range("F11").Value = 'email'
$ws.Range("G11").Value = 'Follow up requested'
# Row 12
$ws.Range("A12").Value = '<2019-03-26 Tue 21:56>'
$ws.Range("B12").Value = 'Systems Integrator VAR'
$ws.Range("C12").Value = '$90K/yr'
$ws.Range("D12").Value = 'http://www.linkedin.com'
$ws.Range("E12").Value = 'David Talamo Financial MGR'
$ws.Range("F12").Value = 'web'
$ws.Range("G12").Value = 'Connected'
# Row 13
$ws.Range("A13").Value = '<2019-03-27 Wed 19:05>'
$ws.Range("B13").Value = 'SQL Server DBA'
$ws.Range("C13").Value = '$90K'
$ws.Range("D13").Value = 'MassHire Framingham Career Centers'
$ws.Range("E13").Value = 'Frank Yeoung'
$ws.Range("F13").Value = 'email'
$ws.Range("G13").Value = 'Contact maintenance and discussing opportunities'
# Row 14
$ws.Range("A14").Value = '<2019-03-27 Wed 20:20>'
$ws.Range("B14").Value = 'Operations/Project Manager with Salesforce experience'
$ws.Range("D14").Value = 'US Tech Solutions.'
$ws.Range("E14").Value = 'Kapil Kumar <kapil@ustechsolutions.com>'
$ws.Range("F14").Value = 'email'
$ws.Range("G14").Value = 'Follow up requested'
# Row 15
$ws.Range("A15").Value = '<2019-03-27 Wed 20:25>'
$ws.Range("B15").Value = 'Deskside Technician'
$ws.Range("D15").Value = 'ApTask'
$ws.Range("E15").Value = 'Dipak Domane dipakd@aptask.com'
$ws.Range("F15").Value = 'email'
$ws.Range("G15").Value = 'Follow up requested'
# Row 16
$ws.Range("A16").Value = '<2019-03-27 Wed 17:58>'
$ws.Range("B16").Value = 'Systems Integrator'
$ws.Range("C16").Value = '$90K/yr'
$ws.Range("D16").Value = 'Tatnuck Worcester Rd Westboro'
$ws.Range("E16").Value = 'David Sullivan and Group tatnuck_group@yahoogroups.com'
$ws.Range("F16").Value = 'email'
$ws.Range("G16").Value = 'Email reminder of Salesforce event'
# Row 17
$ws.Range("A17").Value = '<2019-03-28 Thu 20:51>'
$ws.Range("B17").Value = 'Web Developer'
$ws.Range("D17").Value = 'Talent Retriever'
$ws.Range("E17").Value = 'Leslie Morgan inmail-hit-reply@linkedin.com'
$ws.Range("F17").Value = 'email'
$ws.Range("G17").Value = 'Follow up requested'
# Row 18
$ws.Range("A18").Value = '<2019-03-28 Thu 21:23>'
$ws.Range("B18").Value = 'Helpdesk Temp/Perm'
$ws.Range("D18").Value = 'The CEI Group'
$ws.Range("E18").Value = 'Daniel DesJardins ddesjardins@theceigroup.com'
$ws.Range("F18").Value = 'email'
$ws.Range("G18").Value = 'Follow up requested'
# Row 19
$ws.Range("A19").Value = '<2019-03-28 Thu 23:32>'
$ws.Range("B19").Value = 'Senior Business Data analyst'
$ws.Range("D19").Value = 'IDC Technologies Inc'
$ws.Range("E19").Value = 'Ritesh Jaiswal <ritesh.j@idctechnologies.com>'
$ws.Range("F19").Value = 'email'
$ws.Range("G19").Value = 'Follow up requested'
# Row 20
$ws.Range("A20").Value = '<2019-03-29 Fri 07:46>'
$ws.Range("B20").Value = 'Intermediate Technical Staff'
$ws.Range("D20").Value = 'ZipRecruiter'
$ws.Range("E20").Value = 'Phil @ ZipRecruiter <phil@ziprecruiter.com>'
$ws.Range("F20").Value = 'email'
$ws.Range("G20").Value = 'Follow up requested'
# Row 21
$ws.Range("A21").Value = '<2019-03-29 Fri 07:54>'
$ws.Range("B21").Value = 'PT Evening Instructor'
$ws.Range("D21").Value = 'MassHire'
$ws.Range("E21").Value = 'Ed Lawrence elawrence@masshiresw.com'
$ws.Range("F21").Value = 'email'
$ws.Range("G21").Value = 'Follow up requested'
# Row 22
$ws.Range("A22").Value = '<2019-03-30 Sat 19:57>'
$ws.Range("B22").Value = 'Test Engineer'
$ws.Range("D22").Value = 'E TalentNetwork'
$ws.Range("E22").Value = 'Jeremy Williams <jeremyw@etalentnetwork.com>'
$ws.Range("F22").Value = 'email'
$ws.Range("G22").Value = 'Follow up requested'

# Clear cells that no longer have content (Payrate column in rows now without it)
$ws.Range("C4").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
